$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure value/percent columns keep their original text formatting instead of
# being auto-coerced into numbers by COM when the literal looks numeric.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.294.63'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '2.124.98'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").Value = '347.32'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '0.5222'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").Value = '0.4482'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '54.20'
$ws.Range("E9").Value = '  +4.34%  '
$ws.Range("D10").Value = '0.09388'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").Value = '1.187'
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = '25.49'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '8.693'
$ws.Range("E13").Value = '  +7.81%  '
$ws.Range("D14").Value = '6.979'
$ws.Range("E14").Value = '  +3.44%  '
$ws.Range("D15").Value = '2.121.90'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = '102.97'
$ws.Range("E16").Value = '  +3.73%  '
$ws.Range("D17").Value = '0.00001177'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '1.006'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '21.61'
$ws.Range("E19").Value = '  +5.17%  '
$ws.Range("D20").Value = '0.06706'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '6.321'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").Value = '1.007'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '30.273.04'
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = '12.78'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").Value = '2.338'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").Value = '2.370.53'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '22.26'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").Value = '2.565'
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("D29").Value = '163.48'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '134.47'
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("D31").Value = '1.163'
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").Value = '1.802'
$ws.Range("E32").Value = '  +11.30%  '
$ws.Range("D33").Value = '0.1061'
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").Value = '6.868'
$ws.Range("E34").Value = '  +11.66%  '
$ws.Range("D35").Value = '6.327'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").Value = '3.964'
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").Value = '10.77'
$ws.Range("E37").Value = '  +6.47%  '
$ws.Range("D38").Value = '0.02653'
$ws.Range("E38").Value = '  +3.21%  '
$ws.Range("D39").Value = '0.06882'
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("D40").Value = '0.7176'
$ws.Range("E40").Value = '  +3.74%  '
$ws.Range("D41").Value = '12.82'
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("D42").Value = '0.2259'
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").Value = '1.336'
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("D44").Value = '0.6985'
$ws.Range("E44").Value = '  +4.71%  '
$ws.Range("D45").Value = '14.78'
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("D46").Value = '2.410'
$ws.Range("E46").Value = '  +5.88%  '
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '1.270'
$ws.Range("E48").Value = '  +8.34%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '3.635'
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").Value = '0.00000000348'
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.232'
$ws.Range("E51").Value = '  +1.19%  '
